$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 99: correct the date (E99) and let the shared formula in G99 recalc
$ws.Range("E99").Value = 45754

# Row 100: correct the date (E100) and let the shared formula in G100 recalc
$ws.Range("E100").Value = 45754

# New row 101: a "Tag" action by Hudson Wilson on Connor Scott
$ws.Range("A101").Value = 103
$ws.Range("B101").Value = "Tag"
$ws.Range("C101").Value = "Hudson Wilson"
$ws.Range("D101").Value = "Connor Scott"
$ws.Range("E101").Value = 45754
$ws.Range("F101").Value = 0.45833333333333331

# Match the number formatting/style used by the rest of the Unix column
$ws.Range("G100").Copy()
$ws.Range("G101").PasteSpecial(-4122)
$ws.Range("G101").Formula = "=((E101+F101)-DATE(1970,1,1))*86400"

# Update selection to reflect the newly-added rows
$ws.Range("A97:B101").Select() | Out-Null
